$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 98
$ws_ALC.Range("H98").Value = 570.5
$ws_ALC.Range("I98").Value = 571.4761999999999
$ws_ALC.Range("J98").Value = 550
$ws_ALC.Range("K98").Value = 571.4761999999999
$ws_ALC.Range("L98").Value = 550
$ws_ALC.Range("M98").Value = 926.5238000000001
$ws_ALC.Range("N98").Value = -3546

# ALC row 107
$ws_ALC.Range("H107").Value = 431.73334
$ws_ALC.Range("J107").Value = 496.375
$ws_ALC.Range("L107").Value = 496.375
$ws_ALC.Range("N107").Value = -4336.375

# ALC row 113
$ws_ALC.Range("H113").Value = 34486252
$ws_ALC.Range("I113").Value = 71431650
$ws_ALC.Range("J113").Value = 3883.1333
$ws_ALC.Range("K113").Value = 71431650
$ws_ALC.Range("L113").Value = 3883.1333
$ws_ALC.Range("M113").Value = -71428396
$ws_ALC.Range("N113").Value = -10391.1333

# ALC row 122
$ws_ALC.Range("H122").Value = 570.5
$ws_ALC.Range("I122").Value = 571.4761999999999
$ws_ALC.Range("J122").Value = 550
$ws_ALC.Range("K122").Value = 1714.4286
$ws_ALC.Range("L122").Value = 1650
$ws_ALC.Range("M122").Value = 735.5714000000003
$ws_ALC.Range("N122").Value = -6550

# ALC row 137
$ws_ALC.Range("H137").Value = 1329.7028
$ws_ALC.Range("I137").Value = 1246.4642
$ws_ALC.Range("J137").Value = 1588.6666
$ws_ALC.Range("K137").Value = 3739.3926
$ws_ALC.Range("L137").Value = 4765.9998
$ws_ALC.Range("M137").Value = -1189.3926
$ws_ALC.Range("N137").Value = -9865.9998

# ARM row 2
$ws_ARM.Range("H2").Value = 1712.32
$ws_ARM.Range("I2").Value = 1579.2632
$ws_ARM.Range("K2").Value = 1579.2632
$ws_ARM.Range("M2").Value = -1466.2632

# ARM row 32
$ws_ARM.Range("H32").Value = 6448.9565
$ws_ARM.Range("I32").Value = 4893.9756
$ws_ARM.Range("K32").Value = 4893.9756
$ws_ARM.Range("M32").Value = -4606.9756

# ARM row 116
$ws_ARM.Range("H116").Value = 1712.32
$ws_ARM.Range("I116").Value = 1579.2632
$ws_ARM.Range("K116").Value = 1579.2632
$ws_ARM.Range("M116").Value = 714.7367999999999

# BSM row 3
$ws_BSM.Range("H3").Value = 1712.32
$ws_BSM.Range("I3").Value = 1579.2632
$ws_BSM.Range("K3").Value = 1579.2632
$ws_BSM.Range("M3").Value = -1465.2632

# BSM row 20
$ws_BSM.Range("H20").Value = 5230.5
$ws_BSM.Range("I20").Value = 7323.25
$ws_BSM.Range("J20").Value = 1045
$ws_BSM.Range("K20").Value = 7323.25
$ws_BSM.Range("L20").Value = 1045
$ws_BSM.Range("M20").Value = -7076.25
$ws_BSM.Range("N20").Value = -1539

# CRP row 25
$ws_CRP.Range("H25").Value = 5999.25
$ws_CRP.Range("I25").Value = 6332.3335
$ws_CRP.Range("J25").Value = 5000
$ws_CRP.Range("K25").Value = 6332.3335
$ws_CRP.Range("L25").Value = 5000
$ws_CRP.Range("M25").Value = -6158.3335
$ws_CRP.Range("N25").Value = -5348

# CRP row 31
$ws_CRP.Range("H31").Value = 13703.322
$ws_CRP.Range("I31").Value = 33168.1
$ws_CRP.Range("J31").Value = 4434.381
$ws_CRP.Range("K31").Value = 33168.1
$ws_CRP.Range("L31").Value = 4434.381
$ws_CRP.Range("M31").Value = -32873.1
$ws_CRP.Range("N31").Value = -5024.381

# CRP row 34
$ws_CRP.Range("H34").Value = 13703.322
$ws_CRP.Range("I34").Value = 33168.1
$ws_CRP.Range("J34").Value = 4434.381
$ws_CRP.Range("K34").Value = 33168.1
$ws_CRP.Range("L34").Value = 4434.381
$ws_CRP.Range("M34").Value = -32966.1
$ws_CRP.Range("N34").Value = -4838.381

# CRP row 62
$ws_CRP.Range("H62").Value = 62503130
$ws_CRP.Range("I62").Value = 71431220
$ws_CRP.Range("J62").Value = 6506
$ws_CRP.Range("K62").Value = 71431220
$ws_CRP.Range("L62").Value = 6506
$ws_CRP.Range("M62").Value = -71430596
$ws_CRP.Range("N62").Value = -7754

# CRP row 65
$ws_CRP.Range("H65").Value = 62503130
$ws_CRP.Range("I65").Value = 71431220
$ws_CRP.Range("J65").Value = 6506
$ws_CRP.Range("K65").Value = 357156100
$ws_CRP.Range("L65").Value = 32530
$ws_CRP.Range("M65").Value = -357152980
$ws_CRP.Range("N65").Value = -38770

# CRP row 86
$ws_CRP.Range("H86").Value = 12833956
$ws_CRP.Range("I86").Value = 3587.125
$ws_CRP.Range("J86").Value = 33362548
$ws_CRP.Range("K86").Value = 3587.125
$ws_CRP.Range("L86").Value = 33362548
$ws_CRP.Range("M86").Value = -2464.125
$ws_CRP.Range("N86").Value = -33364794

# CRP row 89
$ws_CRP.Range("H89").Value = 12833956
$ws_CRP.Range("I89").Value = 3587.125
$ws_CRP.Range("J89").Value = 33362548
$ws_CRP.Range("K89").Value = 17935.625
$ws_CRP.Range("L89").Value = 166812740
$ws_CRP.Range("M89").Value = -12319.625
$ws_CRP.Range("N89").Value = -166823972

# CUL row 63
$ws_CUL.Range("H63").Value = 4702.6
$ws_CUL.Range("I63").Value = 1999.5
$ws_CUL.Range("J63").Value = 6504.6665
$ws_CUL.Range("K63").Value = 5998.5
$ws_CUL.Range("L63").Value = 19513.9995
$ws_CUL.Range("M63").Value = -5249.5
$ws_CUL.Range("N63").Value = -21011.9995

# CUL row 66
$ws_CUL.Range("H66").Value = 4702.6
$ws_CUL.Range("I66").Value = 1999.5
$ws_CUL.Range("J66").Value = 6504.6665
$ws_CUL.Range("K66").Value = 17995.5
$ws_CUL.Range("L66").Value = 58541.9985
$ws_CUL.Range("M66").Value = -14251.5
$ws_CUL.Range("N66").Value = -66029.9985

# CUL row 113
$ws_CUL.Range("H113").Value = 11825.556
$ws_CUL.Range("I113").Value = 33930.332
$ws_CUL.Range("J113").Value = 773.1667
$ws_CUL.Range("K113").Value = 101790.996
$ws_CUL.Range("L113").Value = 2319.5001
$ws_CUL.Range("M113").Value = -99620.99600000001
$ws_CUL.Range("N113").Value = -6659.5001

# CUL row 129
$ws_CUL.Range("H129").Value = 1062.2307
$ws_CUL.Range("J129").Value = 1147.9
$ws_CUL.Range("L129").Value = 3443.7
$ws_CUL.Range("N129").Value = -13443.7

# CUL row 131
$ws_CUL.Range("H131").Value = 769.64
$ws_CUL.Range("I131").Value = 275
$ws_CUL.Range("J131").Value = 779.7347
$ws_CUL.Range("K131").Value = 825
$ws_CUL.Range("L131").Value = 2339.2041
$ws_CUL.Range("M131").Value = 4215
$ws_CUL.Range("N131").Value = -12419.2041

# GSM row 39
$ws_GSM.Range("H39").Value = 29999.334
$ws_GSM.Range("J39").Value = 29999.334
$ws_GSM.Range("L39").Value = 29999.334
$ws_GSM.Range("N39").Value = -31063.334

# GSM row 70
$ws_GSM.Range("H70").Value = 4597.4707
$ws_GSM.Range("I70").Value = 4600.8
$ws_GSM.Range("K70").Value = 4600.8
$ws_GSM.Range("M70").Value = -4330.8

# GSM row 73
$ws_GSM.Range("H73").Value = 4597.4707
$ws_GSM.Range("I73").Value = 4600.8
$ws_GSM.Range("K73").Value = 4600.8
$ws_GSM.Range("M73").Value = -3664.8

# LTW row 46
$ws_LTW.Range("H46").Value = 740.9474
$ws_LTW.Range("I46").Value = 859.6
$ws_LTW.Range("J46").Value = 698.5714
$ws_LTW.Range("K46").Value = 859.6
$ws_LTW.Range("L46").Value = 698.5714
$ws_LTW.Range("M46").Value = -671.6
$ws_LTW.Range("N46").Value = -1074.5714

# WVR row 27
$ws_WVR.Range("H27").Value = 35127
$ws_WVR.Range("J27").Value = 35127
$ws_WVR.Range("L27").Value = 35127
$ws_WVR.Range("N27").Value = -35265

# WVR row 100
$ws_WVR.Range("H100").Value = 400.33334
$ws_WVR.Range("I100").Value = 350.5
$ws_WVR.Range("K100").Value = 701
$ws_WVR.Range("M100").Value = -160

# WVR row 133
$ws_WVR.Range("H133").Value = 52715
$ws_WVR.Range("J133").Value = 52715
$ws_WVR.Range("L133").Value = 52715
$ws_WVR.Range("N133").Value = -62835

Write-Host "Applied 167 cell updates across 8 sheets"
